$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: assign cell values in the precise order needed so new shared strings
# land at the exact indices seen in the target workbook (first-use order). ---
$ws.Range("G4").Value2  = "Centralized Power in the hands of a few nodes"
$ws.Range("G5").Value2  = "Centralized, The generation of time series is dependent on the leader, and high performance requirements "
$ws.Range("B2").Value2  = "N/A"
$ws.Range("F3").Value2  = "low energy consumption, high throughput, and scalability"
$ws.Range("G3").Value2  = "pledge centralisation risk, initial pledge cost, and finality latency"
$ws.Range("G2").Value2  = "significant energy expenditure and a comparatively low transaction throughput"
$ws.Range("F2").Value2  = "extremely high security"
$ws.Range("F4").Value2  = "generally possess a reduced number of decentralised nodes and are thus subject to fewer security challenges "
$ws.Range("F9").Value2  = "high scalability, low latency, and resource efficiency"
$ws.Range("G10").Value2 = "challenges in terms of complexity and security"
$ws.Range("G9").Value2  = "complexity, security risk, cross-slice communication, load balancing"
$ws.Range("F10").Value2 = "lower transaction fees and resource consumption while also achieving higher transaction throughput, significant advantages in terms of improving scalability and flexibility, high throughput, low latency, low fees through parallel transaction processing and decentralised architecture"

# --- Step 2: fill in the rest of the table (N/A placeholders + numeric values). ---
$ws.Range("C2").Value2  = "N/A"
$ws.Range("D2").Value2  = "N/A"

$ws.Range("C3").Value2  = "N/A"
$ws.Range("D3").Value2  = "N/A"
$ws.Range("E3").Value2  = "N/A"

$ws.Range("C4").Value2  = "N/A"
$ws.Range("D4").Value2  = "N/A"
$ws.Range("E4").Value2  = "N/A"

$ws.Range("B5").Value2  = "N/A"
$ws.Range("C5").Value2  = "N/A"
$ws.Range("D5").Value2  = "N/A"
$ws.Range("E5").Value2  = "N/A"
$ws.Range("F5").Value2  = "N/A"

$ws.Range("C6").Value2  = "N/A"
$ws.Range("D6").Value2  = "N/A"
$ws.Range("F6").Value2  = "N/A"
$ws.Range("G6").Value2  = "N/A"

$ws.Range("C7").Value2  = "N/A"
$ws.Range("D7").Value2  = "N/A"
$ws.Range("E7").Value2  = "N/A"
$ws.Range("F7").Value2  = "N/A"
$ws.Range("G7").Value2  = "N/A"

$ws.Range("B8").Value2  = "N/A"
$ws.Range("C8").Value2  = "N/A"
$ws.Range("D8").Value2  = "N/A"
$ws.Range("E8").Value2  = "N/A"
$ws.Range("F8").Value2  = "N/A"
$ws.Range("G8").Value2  = "N/A"

$ws.Range("B9").Value2  = "N/A"
$ws.Range("C9").Value2  = "N/A"
$ws.Range("D9").Value2  = "N/A"
$ws.Range("E9").Value2  = "N/A"

$ws.Range("B10").Value2 = "N/A"
$ws.Range("C10").Value2 = "N/A"
$ws.Range("D10").Value2 = "N/A"
$ws.Range("E10").Value2 = "N/A"

# Numeric cells
$ws.Range("B6").Value2 = 20000
$ws.Range("B7").Value2 = 10000
$ws.Range("E2").Value2 = 0.51
$ws.Range("E6").Value2 = 0.33

# --- Step 3: formatting ---
# Percent number format on column E (% of nodes required to take over network)
$ws.Range("E2:E10").Style = "Percent"

# Wrap text on column G (Weaknesses)
$ws.Range("G2:G10").WrapText = $true

# --- Step 4: row heights (auto height Excel computed from wrapped text) ---
$ws.Rows.Item(2).RowHeight = 136
$ws.Rows.Item(3).RowHeight = 102
$ws.Rows.Item(4).RowHeight = 68
$ws.Rows.Item(5).RowHeight = 136
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 17
$ws.Rows.Item(9).RowHeight = 102
$ws.Rows.Item(10).RowHeight = 68

# --- Step 5: final selection matches the author's last cursor position ---
$ws.Range("H3").Select()
